$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4 new values (previously held by row 6)
$ws.Range("D4").Value = 44204
$ws.Range("L4").Value = "Primera"
$ws.Range("M4").Value = 110
$ws.Range("N4").Value = 7000
$ws.Range("O4").Value = 7500
$ws.Range("P4").Value = 7318
$ws.Range("S4").Value = 1045

# Row 5 new values (previously held by row 4)
$ws.Range("L5").Value = "Especial"
$ws.Range("M5").Value = 20
$ws.Range("N5").Value = 15000
$ws.Range("O5").Value = 15000
$ws.Range("P5").Value = 15000
$ws.Range("S5").Value = 2143

# Row 6 new values (previously held by row 5)
$ws.Range("D6").Value = 44189
$ws.Range("M6").Value = 30
$ws.Range("N6").Value = 13000
$ws.Range("O6").Value = 13000
$ws.Range("P6").Value = 13000
$ws.Range("S6").Value = 1857
